{"js": "// The SRS title (\"Subtitle\" paragraph) currently reads:\n//   \"InterTolls - \u03a3\u03cd\u03c3\u03c4\u03b7\u03bc\u03b1 \u03b4\u03b9\u03b1\u03bb\u03b5\u03b9\u03c4\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03cc\u03c4\u03b7\u03c4\u03b1\u03c2 \u03c3\u03c4\u03b1 \u03b4\u03b9\u03cc\u03b4\u03b9\u03b1\"\n// The project / product name changed from \"InterTolls\" to \"Tollways\";\n// the rest of the subtitle stays the same. Only the document body is in\n// scope (the footer also contains the string \"InterTolls \u2013 SoftEng 43\",\n// but that is a different phrase and is not touched by this edit).\n\nconst body = context.document.body;\n\n// Search only the exact, case-sensitive whole word \"InterTolls\" so we\n// don't accidentally match a substring of something else.\nconst results = body.search(\"InterTolls\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"InterTolls\" in the document body.');\n}\n\n// Replace every occurrence found in the body (in practice there is a\n// single one, in the Subtitle paragraph) with the new name.\nfor (const range of results.items) {\n  range.insertText(\"Tollways\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The SRS title (\"Subtitle\" paragraph) currently reads:\n#   \"InterTolls - \u03a3\u03cd\u03c3\u03c4\u03b7\u03bc\u03b1 \u03b4\u03b9\u03b1\u03bb\u03b5\u03b9\u03c4\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03cc\u03c4\u03b7\u03c4\u03b1\u03c2 \u03c3\u03c4\u03b1 \u03b4\u03b9\u03cc\u03b4\u03b9\u03b1\"\n# The product/project name changed from \"InterTolls\" to \"Tollways\"; the\n# rest of the subtitle is unchanged.\n#\n# $d.Content is the main document body story only, so this will not touch\n# the (different) \"InterTolls\" occurrence in the page footer.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"InterTolls\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Tollways\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, `\n    $false, $false, $false, $true, $wdFindContinue, $false, `\n    $find.Replacement.Text, $wdReplaceAll)\n\nif (-not $found) {\n    throw 'Could not find \"InterTolls\" in the document body.'\n}\n"}
